# MyHints.xlsx — add a new "Plus One" problem/hint row (row 9) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A8/B8 = s="3" wrap-text style, C8 = s="2" red-font style) is the
# closest existing row with the exact formatting the new row needs, so copy
# its formats down into row 9 before writing the new values.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("A9:B9").PasteSpecial(-4122) | Out-Null

$ws.Range("C8").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Value = "Plus One"
$ws.Range("B9").Value = "Plus One( Adding 1 to large intDo not attempt to convert from dtring to Int. Instead use for loop from behind/reversed, and for every count in the for loop, use a while loop to pop the last item and add 1. maintain condition if sum > 9"
$ws.Range("C9").Value = "EASY"

# Row 9 needs a taller row (68pt) to fit the wrapped hint text.
$ws.Rows("9").RowHeight = 68

# Keep the previously selected cell as the active selection, and scroll the
# sheet view down so row 4 is pinned at the top (matches topLeftCell="A4").
$ws.Range("C9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
